# "color changes and steel changes"
# Update the base Steel row (row 4) on the Materials sheet: B4:G4 become
# formulas derived from row 2 (same pattern already used elsewhere in the
# sheet, e.g. row 21 uses (row3*2)*0.75). Everything else in the workbook
# (Part 1 / Part 2 sheets, and the other summary rows on Materials) already
# references these cells via formulas, so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$ws.Range("B4").Formula = "=(B2*2)*0.75"
$ws.Range("C4").Formula = "=(C2*2)*0.75"
$ws.Range("D4").Formula = "=(D2*2)*0.75"
$ws.Range("E4").Formula = "=(E2*2)*0.75"
$ws.Range("F4").Formula = "=(F2*2)*0.75"
$ws.Range("G4").Formula = "=(G2*2)*0.75"

# Match the "Invoer" style used elsewhere (e.g. row 16) that shows two
# decimal places, since the new computed values are no longer whole numbers.
$ws.Range("B4:G4").NumberFormat = "0.00"

# Leave the selection where the author last left it when saving.
$ws.Range("M8").Select()
